$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: re-point existing rows (144-242) to their corrected match data ---
# Column A (row id) stays put; B..AC are rewritten in place per row.
# (C/D are always the constant "Germany Regionalliga West" label, so they are skipped.)

# row 146
$ws.Cells.Item(146, 2).Value = 6886952
$ws.Cells.Item(146, 5).Value = 45255.41666666666
$ws.Cells.Item(146, 6).Value = "RotWeiss Oberhausen"
$ws.Cells.Item(146, 7).Value = "SC Paderborn 07 II"
$ws.Cells.Item(146, 8).Value = 4
$ws.Cells.Item(146, 9).Value = 1
$ws.Cells.Item(146, 10).Value = "H"
$ws.Cells.Item(146, 11).Value = 2
$ws.Cells.Item(146, 12).Value = 3.75
$ws.Cells.Item(146, 13).Value = 2.875
$ws.Cells.Item(146, 14).Value = 1.615
$ws.Cells.Item(146, 15).Value = 4.2
$ws.Cells.Item(146, 16).Value = 4
$ws.Cells.Item(146, 17).Value = -0.75
$ws.Cells.Item(146, 18).Value = 1.775
$ws.Cells.Item(146, 19).Value = 2.025
$ws.Cells.Item(146, 20).Value = 3.25
$ws.Cells.Item(146, 21).Value = 2
$ws.Cells.Item(146, 22).Value = 1.8
$ws.Cells.Item(146, 23).Value = 0.615
$ws.Cells.Item(146, 24).Value = -1
$ws.Cells.Item(146, 25).Value = -1
$ws.Cells.Item(146, 26).Value = 0.7749999999999999
$ws.Cells.Item(146, 27).Value = -1
$ws.Cells.Item(146, 28).Value = 1
$ws.Cells.Item(146, 29).Value = -1

# row 149
$ws.Cells.Item(149, 2).Value = 6886950
$ws.Cells.Item(149, 5).Value = 45255.41666666666
$ws.Cells.Item(149, 6).Value = "SC Fortuna Kln"
$ws.Cells.Item(149, 7).Value = "SV Rodinghausen"
$ws.Cells.Item(149, 8).Value = 1
$ws.Cells.Item(149, 9).Value = 3
$ws.Cells.Item(149, 10).Value = "A"
$ws.Cells.Item(149, 11).Value = 1.363
$ws.Cells.Item(149, 12).Value = 4.75
$ws.Cells.Item(149, 13).Value = 6
$ws.Cells.Item(149, 14).Value = 1.75
$ws.Cells.Item(149, 15).Value = 3.6
$ws.Cells.Item(149, 16).Value = 3.75
$ws.Cells.Item(149, 17).Value = -0.5
$ws.Cells.Item(149, 18).Value = 1.8
$ws.Cells.Item(149, 19).Value = 2
$ws.Cells.Item(149, 20).Value = 2.75
$ws.Cells.Item(149, 21).Value = 2
$ws.Cells.Item(149, 22).Value = 1.8
$ws.Cells.Item(149, 23).Value = -1
$ws.Cells.Item(149, 24).Value = -1
$ws.Cells.Item(149, 25).Value = 2.75
$ws.Cells.Item(149, 26).Value = -1
$ws.Cells.Item(149, 27).Value = 1
$ws.Cells.Item(149, 28).Value = 1
$ws.Cells.Item(149, 29).Value = -1

# row 167
$ws.Cells.Item(167, 2).Value = 6880544
$ws.Cells.Item(167, 5).Value = 45325.41666666666
$ws.Cells.Item(167, 6).Value = "1 FC Bocholt"
$ws.Cells.Item(167, 7).Value = "Cologne II"
$ws.Cells.Item(167, 8).Value = 2
$ws.Cells.Item(167, 9).Value = 2
$ws.Cells.Item(167, 10).Value = "D"
$ws.Cells.Item(167, 11).Value = 2
$ws.Cells.Item(167, 12).Value = 3.6
$ws.Cells.Item(167, 13).Value = 3
$ws.Cells.Item(167, 14).Value = 2.1
$ws.Cells.Item(167, 15).Value = 3.6
$ws.Cells.Item(167, 16).Value = 2.8
$ws.Cells.Item(167, 17).Value = -0.25
$ws.Cells.Item(167, 18).Value = 1.925
$ws.Cells.Item(167, 19).Value = 1.925
$ws.Cells.Item(167, 20).Value = 2.75
$ws.Cells.Item(167, 21).Value = 1.8
$ws.Cells.Item(167, 22).Value = 2.05
$ws.Cells.Item(167, 23).Value = -1
$ws.Cells.Item(167, 24).Value = 2.6
$ws.Cells.Item(167, 25).Value = -1
$ws.Cells.Item(167, 26).Value = -0.5
$ws.Cells.Item(167, 27).Value = 0.4625
$ws.Cells.Item(167, 28).Value = 0.8
$ws.Cells.Item(167, 29).Value = -1

# row 168
$ws.Cells.Item(168, 2).Value = 6882764
$ws.Cells.Item(168, 5).Value = 45325.41666666666
$ws.Cells.Item(168, 6).Value = "Gutersloh 2000"
$ws.Cells.Item(168, 7).Value = "SSVg Velbert"
$ws.Cells.Item(168, 8).Value = 1
$ws.Cells.Item(168, 9).Value = 2
$ws.Cells.Item(168, 10).Value = "A"
$ws.Cells.Item(168, 11).Value = 1.666
$ws.Cells.Item(168, 12).Value = 3.75
$ws.Cells.Item(168, 13).Value = 4
$ws.Cells.Item(168, 14).Value = 1.833
$ws.Cells.Item(168, 15).Value = 3.6
$ws.Cells.Item(168, 16).Value = 3.6
$ws.Cells.Item(168, 17).Value = -0.5
$ws.Cells.Item(168, 18).Value = 1.875
$ws.Cells.Item(168, 19).Value = 1.975
$ws.Cells.Item(168, 20).Value = 2.75
$ws.Cells.Item(168, 21).Value = 1.825
$ws.Cells.Item(168, 22).Value = 2.025
$ws.Cells.Item(168, 23).Value = -1
$ws.Cells.Item(168, 24).Value = -1
$ws.Cells.Item(168, 25).Value = 2.6
$ws.Cells.Item(168, 26).Value = -1
$ws.Cells.Item(168, 27).Value = 0.9750000000000001
$ws.Cells.Item(168, 28).Value = 0.4125
$ws.Cells.Item(168, 29).Value = -0.5

# row 169
$ws.Cells.Item(169, 2).Value = 6884503
$ws.Cells.Item(169, 5).Value = 45325.41666666666
$ws.Cells.Item(169, 6).Value = "Borussia Mgladbach II"
$ws.Cells.Item(169, 7).Value = "Alemannia Aachen"
$ws.Cells.Item(169, 8).Value = 0
$ws.Cells.Item(169, 9).Value = 4
$ws.Cells.Item(169, 10).Value = "A"
$ws.Cells.Item(169, 11).Value = 3.1
$ws.Cells.Item(169, 12).Value = 3.75
$ws.Cells.Item(169, 13).Value = 1.909
$ws.Cells.Item(169, 14).Value = 3.75
$ws.Cells.Item(169, 15).Value = 4.2
$ws.Cells.Item(169, 16).Value = 1.666
$ws.Cells.Item(169, 17).Value = 0.75
$ws.Cells.Item(169, 18).Value = 1.9
$ws.Cells.Item(169, 19).Value = 1.9
$ws.Cells.Item(169, 20).Value = 2.75
$ws.Cells.Item(169, 21).Value = 1.775
$ws.Cells.Item(169, 22).Value = 2.025
$ws.Cells.Item(169, 23).Value = -1
$ws.Cells.Item(169, 24).Value = -1
$ws.Cells.Item(169, 25).Value = 0.6659999999999999
$ws.Cells.Item(169, 26).Value = -1
$ws.Cells.Item(169, 27).Value = 0.8999999999999999
$ws.Cells.Item(169, 28).Value = 0.7749999999999999
$ws.Cells.Item(169, 29).Value = -1

# row 170
$ws.Cells.Item(170, 2).Value = 6884536
$ws.Cells.Item(170, 5).Value = 45325.41666666666
$ws.Cells.Item(170, 6).Value = "Wuppertaler"
$ws.Cells.Item(170, 7).Value = "Schalke II"
$ws.Cells.Item(170, 8).Value = 3
$ws.Cells.Item(170, 9).Value = 1
$ws.Cells.Item(170, 10).Value = "H"
$ws.Cells.Item(170, 11).Value = 1.666
$ws.Cells.Item(170, 12).Value = 3.75
$ws.Cells.Item(170, 13).Value = 4
$ws.Cells.Item(170, 14).Value = 1.85
$ws.Cells.Item(170, 15).Value = 3.6
$ws.Cells.Item(170, 16).Value = 3.6
$ws.Cells.Item(170, 17).Value = -0.5
$ws.Cells.Item(170, 18).Value = 1.9
$ws.Cells.Item(170, 19).Value = 1.9
$ws.Cells.Item(170, 20).Value = 3
$ws.Cells.Item(170, 21).Value = 1.975
$ws.Cells.Item(170, 22).Value = 1.825
$ws.Cells.Item(170, 23).Value = 0.8500000000000001
$ws.Cells.Item(170, 24).Value = -1
$ws.Cells.Item(170, 25).Value = -1
$ws.Cells.Item(170, 26).Value = 0.8999999999999999
$ws.Cells.Item(170, 27).Value = -1
$ws.Cells.Item(170, 28).Value = 0.9750000000000001
$ws.Cells.Item(170, 29).Value = -1

# row 172
$ws.Cells.Item(172, 2).Value = 7764092
$ws.Cells.Item(172, 5).Value = 45328.64583333334
$ws.Cells.Item(172, 6).Value = "SSVg Velbert"
$ws.Cells.Item(172, 7).Value = "RotWeiss Oberhausen"
$ws.Cells.Item(172, 8).Value = 0
$ws.Cells.Item(172, 9).Value = 3
$ws.Cells.Item(172, 10).Value = "A"
$ws.Cells.Item(172, 11).Value = 4.75
$ws.Cells.Item(172, 12).Value = 4
$ws.Cells.Item(172, 13).Value = 1.571
$ws.Cells.Item(172, 14).Value = 5
$ws.Cells.Item(172, 15).Value = 4
$ws.Cells.Item(172, 16).Value = 1.571
$ws.Cells.Item(172, 17).Value = 1
$ws.Cells.Item(172, 18).Value = 1.825
$ws.Cells.Item(172, 19).Value = 1.975
$ws.Cells.Item(172, 20).Value = 2.75
$ws.Cells.Item(172, 21).Value = 1.775
$ws.Cells.Item(172, 22).Value = 2.025
$ws.Cells.Item(172, 23).Value = -1
$ws.Cells.Item(172, 24).Value = -1
$ws.Cells.Item(172, 25).Value = 0.571
$ws.Cells.Item(172, 26).Value = -1
$ws.Cells.Item(172, 27).Value = 0.9750000000000001
$ws.Cells.Item(172, 28).Value = 0.3875
$ws.Cells.Item(172, 29).Value = -0.5

# row 173
$ws.Cells.Item(173, 2).Value = 7691489
$ws.Cells.Item(173, 5).Value = 45328.64583333334
$ws.Cells.Item(173, 6).Value = "SC Wiedenbruck"
$ws.Cells.Item(173, 7).Value = "Gutersloh 2000"
$ws.Cells.Item(173, 8).Value = 0
$ws.Cells.Item(173, 9).Value = 2
$ws.Cells.Item(173, 10).Value = "A"
$ws.Cells.Item(173, 11).Value = 2.15
$ws.Cells.Item(173, 12).Value = 3.25
$ws.Cells.Item(173, 13).Value = 2.9
$ws.Cells.Item(173, 14).Value = 2.15
$ws.Cells.Item(173, 15).Value = 3.4
$ws.Cells.Item(173, 16).Value = 3
$ws.Cells.Item(173, 17).Value = -0.25
$ws.Cells.Item(173, 18).Value = 1.925
$ws.Cells.Item(173, 19).Value = 1.875
$ws.Cells.Item(173, 20).Value = 2.75
$ws.Cells.Item(173, 21).Value = 1.975
$ws.Cells.Item(173, 22).Value = 1.825
$ws.Cells.Item(173, 23).Value = -1
$ws.Cells.Item(173, 24).Value = -1
$ws.Cells.Item(173, 25).Value = 2
$ws.Cells.Item(173, 26).Value = -1
$ws.Cells.Item(173, 27).Value = 0.875
$ws.Cells.Item(173, 28).Value = -1
$ws.Cells.Item(173, 29).Value = 0.825

# row 176
$ws.Cells.Item(176, 2).Value = 6884537
$ws.Cells.Item(176, 5).Value = 45332.41666666666
$ws.Cells.Item(176, 6).Value = "SV Rodinghausen"
$ws.Cells.Item(176, 7).Value = "Wuppertaler"
$ws.Cells.Item(176, 8).Value = 0
$ws.Cells.Item(176, 9).Value = 4
$ws.Cells.Item(176, 10).Value = "A"
$ws.Cells.Item(176, 11).Value = 2.2
$ws.Cells.Item(176, 12).Value = 3.4
$ws.Cells.Item(176, 13).Value = 1.159
$ws.Cells.Item(176, 14).Value = 2.3
$ws.Cells.Item(176, 15).Value = 3.3
$ws.Cells.Item(176, 16).Value = 2.75
$ws.Cells.Item(176, 17).Value = -0.25
$ws.Cells.Item(176, 18).Value = 2.1
$ws.Cells.Item(176, 19).Value = 1.775
$ws.Cells.Item(176, 20).Value = 2.75
$ws.Cells.Item(176, 21).Value = 2.025
$ws.Cells.Item(176, 22).Value = 1.825
$ws.Cells.Item(176, 23).Value = -1
$ws.Cells.Item(176, 24).Value = -1
$ws.Cells.Item(176, 25).Value = 1.75
$ws.Cells.Item(176, 26).Value = -1
$ws.Cells.Item(176, 27).Value = 0.7749999999999999
$ws.Cells.Item(176, 28).Value = 1.025
$ws.Cells.Item(176, 29).Value = -1

# row 177
$ws.Cells.Item(177, 2).Value = 6884504
$ws.Cells.Item(177, 5).Value = 45332.41666666666
$ws.Cells.Item(177, 6).Value = "Alemannia Aachen"
$ws.Cells.Item(177, 7).Value = "RotWeiss Oberhausen"
$ws.Cells.Item(177, 8).Value = 3
$ws.Cells.Item(177, 9).Value = 1
$ws.Cells.Item(177, 10).Value = "H"
$ws.Cells.Item(177, 11).Value = 1.727
$ws.Cells.Item(177, 12).Value = 4
$ws.Cells.Item(177, 13).Value = 3.5
$ws.Cells.Item(177, 14).Value = 2.05
$ws.Cells.Item(177, 15).Value = 3.8
$ws.Cells.Item(177, 16).Value = 2.75
$ws.Cells.Item(177, 17).Value = -0.25
$ws.Cells.Item(177, 18).Value = 1.875
$ws.Cells.Item(177, 19).Value = 1.975
$ws.Cells.Item(177, 20).Value = 2.75
$ws.Cells.Item(177, 21).Value = 2
$ws.Cells.Item(177, 22).Value = 1.85
$ws.Cells.Item(177, 23).Value = 1.05
$ws.Cells.Item(177, 24).Value = -1
$ws.Cells.Item(177, 25).Value = -1
$ws.Cells.Item(177, 26).Value = 0.875
$ws.Cells.Item(177, 27).Value = -1
$ws.Cells.Item(177, 28).Value = 1
$ws.Cells.Item(177, 29).Value = -1

# row 178
$ws.Cells.Item(178, 2).Value = 6886968
$ws.Cells.Item(178, 5).Value = 45332.41666666666
$ws.Cells.Item(178, 6).Value = "SC Paderborn 07 II"
$ws.Cells.Item(178, 7).Value = "Rot Weiss Ahlen"
$ws.Cells.Item(178, 8).Value = 2
$ws.Cells.Item(178, 9).Value = 2
$ws.Cells.Item(178, 10).Value = "D"
$ws.Cells.Item(178, 11).Value = 1.615
$ws.Cells.Item(178, 12).Value = 4
$ws.Cells.Item(178, 13).Value = 4
$ws.Cells.Item(178, 14).Value = 1.85
$ws.Cells.Item(178, 15).Value = 3.75
$ws.Cells.Item(178, 16).Value = 3.2
$ws.Cells.Item(178, 17).Value = -0.5
$ws.Cells.Item(178, 18).Value = 2.05
$ws.Cells.Item(178, 19).Value = 1.8
$ws.Cells.Item(178, 20).Value = 3
$ws.Cells.Item(178, 21).Value = 1.875
$ws.Cells.Item(178, 22).Value = 1.975
$ws.Cells.Item(178, 23).Value = -1
$ws.Cells.Item(178, 24).Value = 2.75
$ws.Cells.Item(178, 25).Value = -1
$ws.Cells.Item(178, 26).Value = -1
$ws.Cells.Item(178, 27).Value = 0.8
$ws.Cells.Item(178, 28).Value = 0.875
$ws.Cells.Item(178, 29).Value = -1

# row 190
$ws.Cells.Item(190, 2).Value = 6885378
$ws.Cells.Item(190, 5).Value = 45339.41666666666
$ws.Cells.Item(190, 6).Value = "SC Fortuna Kln"
$ws.Cells.Item(190, 7).Value = "SSVg Velbert"
$ws.Cells.Item(190, 8).Value = 2
$ws.Cells.Item(190, 9).Value = 1
$ws.Cells.Item(190, 10).Value = "H"
$ws.Cells.Item(190, 11).Value = 1.222
$ws.Cells.Item(190, 12).Value = 5.5
$ws.Cells.Item(190, 13).Value = 9
$ws.Cells.Item(190, 14).Value = 1.333
$ws.Cells.Item(190, 15).Value = 5
$ws.Cells.Item(190, 16).Value = 6.5
$ws.Cells.Item(190, 17).Value = -1.5
$ws.Cells.Item(190, 18).Value = 1.925
$ws.Cells.Item(190, 19).Value = 1.875
$ws.Cells.Item(190, 20).Value = 3
$ws.Cells.Item(190, 21).Value = 1.775
$ws.Cells.Item(190, 22).Value = 2.025
$ws.Cells.Item(190, 23).Value = 0.333
$ws.Cells.Item(190, 24).Value = -1
$ws.Cells.Item(190, 25).Value = -1
$ws.Cells.Item(190, 26).Value = -1
$ws.Cells.Item(190, 27).Value = 0.875
$ws.Cells.Item(190, 28).Value = 0
$ws.Cells.Item(190, 29).Value = -0

# row 191
$ws.Cells.Item(191, 2).Value = 6886972
$ws.Cells.Item(191, 5).Value = 45339.41666666666
$ws.Cells.Item(191, 6).Value = "Rot Weiss Ahlen"
$ws.Cells.Item(191, 7).Value = "SC Wiedenbruck"
$ws.Cells.Item(191, 8).Value = 0
$ws.Cells.Item(191, 9).Value = 1
$ws.Cells.Item(191, 10).Value = "A"
$ws.Cells.Item(191, 11).Value = 2.55
$ws.Cells.Item(191, 12).Value = 3.75
$ws.Cells.Item(191, 13).Value = 2.2
$ws.Cells.Item(191, 14).Value = 2.3
$ws.Cells.Item(191, 15).Value = 3.6
$ws.Cells.Item(191, 16).Value = 2.5
$ws.Cells.Item(191, 17).Value = 0
$ws.Cells.Item(191, 18).Value = 1.825
$ws.Cells.Item(191, 19).Value = 1.975
$ws.Cells.Item(191, 20).Value = 2.5
$ws.Cells.Item(191, 21).Value = 1.825
$ws.Cells.Item(191, 22).Value = 1.975
$ws.Cells.Item(191, 23).Value = -1
$ws.Cells.Item(191, 24).Value = -1
$ws.Cells.Item(191, 25).Value = 1.5
$ws.Cells.Item(191, 26).Value = -1
$ws.Cells.Item(191, 27).Value = 0.9750000000000001
$ws.Cells.Item(191, 28).Value = -1
$ws.Cells.Item(191, 29).Value = 0.9750000000000001

# row 192
$ws.Cells.Item(192, 2).Value = 6882762
$ws.Cells.Item(192, 5).Value = 45339.41666666666
$ws.Cells.Item(192, 6).Value = "1 FC Bocholt"
$ws.Cells.Item(192, 7).Value = "SC Paderborn 07 II"
$ws.Cells.Item(192, 8).Value = 0
$ws.Cells.Item(192, 9).Value = 1
$ws.Cells.Item(192, 10).Value = "A"
$ws.Cells.Item(192, 11).Value = 1.615
$ws.Cells.Item(192, 12).Value = 4
$ws.Cells.Item(192, 13).Value = 4
$ws.Cells.Item(192, 14).Value = 1.85
$ws.Cells.Item(192, 15).Value = 3.6
$ws.Cells.Item(192, 16).Value = 3.25
$ws.Cells.Item(192, 17).Value = -0.5
$ws.Cells.Item(192, 18).Value = 1.95
$ws.Cells.Item(192, 19).Value = 1.85
$ws.Cells.Item(192, 20).Value = 2.75
$ws.Cells.Item(192, 21).Value = 1.975
$ws.Cells.Item(192, 22).Value = 1.825
$ws.Cells.Item(192, 23).Value = -1
$ws.Cells.Item(192, 24).Value = -1
$ws.Cells.Item(192, 25).Value = 2.25
$ws.Cells.Item(192, 26).Value = -1
$ws.Cells.Item(192, 27).Value = 0.8500000000000001
$ws.Cells.Item(192, 28).Value = -1
$ws.Cells.Item(192, 29).Value = 0.825

# row 194
$ws.Cells.Item(194, 2).Value = 6886980
$ws.Cells.Item(194, 5).Value = 45345.64583333334
$ws.Cells.Item(194, 6).Value = "SC Paderborn 07 II"
$ws.Cells.Item(194, 7).Value = "Duren"
$ws.Cells.Item(194, 8).Value = 2
$ws.Cells.Item(194, 9).Value = 0
$ws.Cells.Item(194, 10).Value = "H"
$ws.Cells.Item(194, 11).Value = 1.909
$ws.Cells.Item(194, 12).Value = 3.6
$ws.Cells.Item(194, 13).Value = 3.2
$ws.Cells.Item(194, 14).Value = 2.375
$ws.Cells.Item(194, 15).Value = 3.4
$ws.Cells.Item(194, 16).Value = 2.6
$ws.Cells.Item(194, 17).Value = 0
$ws.Cells.Item(194, 18).Value = 1.825
$ws.Cells.Item(194, 19).Value = 2.025
$ws.Cells.Item(194, 20).Value = 2.75
$ws.Cells.Item(194, 21).Value = 1.85
$ws.Cells.Item(194, 22).Value = 2
$ws.Cells.Item(194, 23).Value = 1.375
$ws.Cells.Item(194, 24).Value = -1
$ws.Cells.Item(194, 25).Value = -1
$ws.Cells.Item(194, 26).Value = 0.825
$ws.Cells.Item(194, 27).Value = -1
$ws.Cells.Item(194, 28).Value = -1
$ws.Cells.Item(194, 29).Value = 1

# row 195
$ws.Cells.Item(195, 2).Value = 6885379
$ws.Cells.Item(195, 5).Value = 45345.64583333334
$ws.Cells.Item(195, 6).Value = "SSVg Velbert"
$ws.Cells.Item(195, 7).Value = "Rot Weiss Ahlen"
$ws.Cells.Item(195, 8).Value = 3
$ws.Cells.Item(195, 9).Value = 2
$ws.Cells.Item(195, 10).Value = "H"
$ws.Cells.Item(195, 11).Value = 2.2
$ws.Cells.Item(195, 12).Value = 3.4
$ws.Cells.Item(195, 13).Value = 2.75
$ws.Cells.Item(195, 14).Value = 2.875
$ws.Cells.Item(195, 15).Value = 3.25
$ws.Cells.Item(195, 16).Value = 2.25
$ws.Cells.Item(195, 17).Value = 0.25
$ws.Cells.Item(195, 18).Value = 1.8
$ws.Cells.Item(195, 19).Value = 2
$ws.Cells.Item(195, 20).Value = 2.75
$ws.Cells.Item(195, 21).Value = 1.85
$ws.Cells.Item(195, 22).Value = 1.95
$ws.Cells.Item(195, 23).Value = 1.875
$ws.Cells.Item(195, 24).Value = -1
$ws.Cells.Item(195, 25).Value = -1
$ws.Cells.Item(195, 26).Value = 0.8
$ws.Cells.Item(195, 27).Value = -1
$ws.Cells.Item(195, 28).Value = 0.8500000000000001
$ws.Cells.Item(195, 29).Value = -1

# row 196
$ws.Cells.Item(196, 2).Value = 6880546
$ws.Cells.Item(196, 5).Value = 45346.41666666666
$ws.Cells.Item(196, 6).Value = "SC Wiedenbruck"
$ws.Cells.Item(196, 7).Value = "1 FC Bocholt"
$ws.Cells.Item(196, 8).Value = 0
$ws.Cells.Item(196, 9).Value = 0
$ws.Cells.Item(196, 10).Value = "D"
$ws.Cells.Item(196, 11).Value = 2.5
$ws.Cells.Item(196, 12).Value = 3.6
$ws.Cells.Item(196, 13).Value = 2.3
$ws.Cells.Item(196, 14).Value = 2.55
$ws.Cells.Item(196, 15).Value = 3.4
$ws.Cells.Item(196, 16).Value = 2.375
$ws.Cells.Item(196, 17).Value = 0
$ws.Cells.Item(196, 18).Value = 1.975
$ws.Cells.Item(196, 19).Value = 1.825
$ws.Cells.Item(196, 20).Value = 2.25
$ws.Cells.Item(196, 21).Value = 1.8
$ws.Cells.Item(196, 22).Value = 2
$ws.Cells.Item(196, 23).Value = -1
$ws.Cells.Item(196, 24).Value = 2.4
$ws.Cells.Item(196, 25).Value = -1
$ws.Cells.Item(196, 26).Value = 0
$ws.Cells.Item(196, 27).Value = -0
$ws.Cells.Item(196, 28).Value = -1
$ws.Cells.Item(196, 29).Value = 1

# row 197
$ws.Cells.Item(197, 2).Value = 6882760
$ws.Cells.Item(197, 5).Value = 45346.41666666666
$ws.Cells.Item(197, 6).Value = "Alemannia Aachen"
$ws.Cells.Item(197, 7).Value = "Gutersloh 2000"
$ws.Cells.Item(197, 8).Value = 4
$ws.Cells.Item(197, 9).Value = 0
$ws.Cells.Item(197, 10).Value = "H"
$ws.Cells.Item(197, 11).Value = 1.4
$ws.Cells.Item(197, 12).Value = 4.5
$ws.Cells.Item(197, 13).Value = 5.5
$ws.Cells.Item(197, 14).Value = 1.4
$ws.Cells.Item(197, 15).Value = 4.333
$ws.Cells.Item(197, 16).Value = 6.5
$ws.Cells.Item(197, 17).Value = -1.25
$ws.Cells.Item(197, 18).Value = 1.825
$ws.Cells.Item(197, 19).Value = 1.975
$ws.Cells.Item(197, 20).Value = 3
$ws.Cells.Item(197, 21).Value = 1.875
$ws.Cells.Item(197, 22).Value = 1.925
$ws.Cells.Item(197, 23).Value = 0.3999999999999999
$ws.Cells.Item(197, 24).Value = -1
$ws.Cells.Item(197, 25).Value = -1
$ws.Cells.Item(197, 26).Value = 0.825
$ws.Cells.Item(197, 27).Value = -1
$ws.Cells.Item(197, 28).Value = 0.875
$ws.Cells.Item(197, 29).Value = -1

# row 198
$ws.Cells.Item(198, 2).Value = 6886978
$ws.Cells.Item(198, 5).Value = 45346.41666666666
$ws.Cells.Item(198, 6).Value = "SV Rodinghausen"
$ws.Cells.Item(198, 7).Value = "Schalke II"
$ws.Cells.Item(198, 8).Value = 0
$ws.Cells.Item(198, 9).Value = 0
$ws.Cells.Item(198, 10).Value = "D"
$ws.Cells.Item(198, 11).Value = 1.727
$ws.Cells.Item(198, 12).Value = 3.5
$ws.Cells.Item(198, 13).Value = 4
$ws.Cells.Item(198, 14).Value = 2.2
$ws.Cells.Item(198, 15).Value = 3.2
$ws.Cells.Item(198, 16).Value = 2.9
$ws.Cells.Item(198, 17).Value = -0.25
$ws.Cells.Item(198, 18).Value = 2
$ws.Cells.Item(198, 19).Value = 1.85
$ws.Cells.Item(198, 20).Value = 2.75
$ws.Cells.Item(198, 21).Value = 1.975
$ws.Cells.Item(198, 22).Value = 1.875
$ws.Cells.Item(198, 23).Value = -1
$ws.Cells.Item(198, 24).Value = 2.2
$ws.Cells.Item(198, 25).Value = -1
$ws.Cells.Item(198, 26).Value = -0.5
$ws.Cells.Item(198, 27).Value = 0.425
$ws.Cells.Item(198, 28).Value = -1
$ws.Cells.Item(198, 29).Value = 0.875

# row 199
$ws.Cells.Item(199, 2).Value = 6886976
$ws.Cells.Item(199, 5).Value = 45346.41666666666
$ws.Cells.Item(199, 6).Value = "SV Lippstadt 08"
$ws.Cells.Item(199, 7).Value = "SC Fortuna Kln"
$ws.Cells.Item(199, 8).Value = 2
$ws.Cells.Item(199, 9).Value = 2
$ws.Cells.Item(199, 10).Value = "D"
$ws.Cells.Item(199, 11).Value = 4
$ws.Cells.Item(199, 12).Value = 3.5
$ws.Cells.Item(199, 13).Value = 1.727
$ws.Cells.Item(199, 14).Value = 5.25
$ws.Cells.Item(199, 15).Value = 3.5
$ws.Cells.Item(199, 16).Value = 1.571
$ws.Cells.Item(199, 17).Value = 0.75
$ws.Cells.Item(199, 18).Value = 1.9
$ws.Cells.Item(199, 19).Value = 1.95
$ws.Cells.Item(199, 20).Value = 2.75
$ws.Cells.Item(199, 21).Value = 1.925
$ws.Cells.Item(199, 22).Value = 1.925
$ws.Cells.Item(199, 23).Value = -1
$ws.Cells.Item(199, 24).Value = 2.5
$ws.Cells.Item(199, 25).Value = -1
$ws.Cells.Item(199, 26).Value = 0.8999999999999999
$ws.Cells.Item(199, 27).Value = -1
$ws.Cells.Item(199, 28).Value = 0.925
$ws.Cells.Item(199, 29).Value = -1

# row 204
$ws.Cells.Item(204, 2).Value = 6886982
$ws.Cells.Item(204, 5).Value = 45353.41666666666
$ws.Cells.Item(204, 6).Value = "Borussia Mgladbach II"
$ws.Cells.Item(204, 7).Value = "FC WegbergBeeck"
$ws.Cells.Item(204, 8).Value = 2
$ws.Cells.Item(204, 9).Value = 0
$ws.Cells.Item(204, 10).Value = "H"
$ws.Cells.Item(204, 11).Value = 1.727
$ws.Cells.Item(204, 12).Value = 3.75
$ws.Cells.Item(204, 13).Value = 3.75
$ws.Cells.Item(204, 14).Value = 1.55
$ws.Cells.Item(204, 15).Value = 4
$ws.Cells.Item(204, 16).Value = 5
$ws.Cells.Item(204, 17).Value = -1
$ws.Cells.Item(204, 18).Value = 1.95
$ws.Cells.Item(204, 19).Value = 1.85
$ws.Cells.Item(204, 20).Value = 2.75
$ws.Cells.Item(204, 21).Value = 1.8
$ws.Cells.Item(204, 22).Value = 2
$ws.Cells.Item(204, 23).Value = 0.55
$ws.Cells.Item(204, 24).Value = -1
$ws.Cells.Item(204, 25).Value = -1
$ws.Cells.Item(204, 26).Value = 0.95
$ws.Cells.Item(204, 27).Value = -1
$ws.Cells.Item(204, 28).Value = -1
$ws.Cells.Item(204, 29).Value = 1

# row 205
$ws.Cells.Item(205, 2).Value = 6886985
$ws.Cells.Item(205, 5).Value = 45353.41666666666
$ws.Cells.Item(205, 6).Value = "Rot Weiss Ahlen"
$ws.Cells.Item(205, 7).Value = "SV Lippstadt 08"
$ws.Cells.Item(205, 8).Value = 1
$ws.Cells.Item(205, 9).Value = 2
$ws.Cells.Item(205, 10).Value = "A"
$ws.Cells.Item(205, 11).Value = 1.85
$ws.Cells.Item(205, 12).Value = 3.75
$ws.Cells.Item(205, 13).Value = 3.3
$ws.Cells.Item(205, 14).Value = 1.909
$ws.Cells.Item(205, 15).Value = 3.75
$ws.Cells.Item(205, 16).Value = 3.3
$ws.Cells.Item(205, 17).Value = -0.5
$ws.Cells.Item(205, 18).Value = 1.975
$ws.Cells.Item(205, 19).Value = 1.875
$ws.Cells.Item(205, 20).Value = 2.75
$ws.Cells.Item(205, 21).Value = 1.825
$ws.Cells.Item(205, 22).Value = 2.025
$ws.Cells.Item(205, 23).Value = -1
$ws.Cells.Item(205, 24).Value = -1
$ws.Cells.Item(205, 25).Value = 2.3
$ws.Cells.Item(205, 26).Value = -1
$ws.Cells.Item(205, 27).Value = 0.875
$ws.Cells.Item(205, 28).Value = 0.4125
$ws.Cells.Item(205, 29).Value = -0.5

# row 206
$ws.Cells.Item(206, 2).Value = 6886984
$ws.Cells.Item(206, 5).Value = 45353.41666666666
$ws.Cells.Item(206, 6).Value = "SC Fortuna Kln"
$ws.Cells.Item(206, 7).Value = "RotWeiss Oberhausen"
$ws.Cells.Item(206, 8).Value = 2
$ws.Cells.Item(206, 9).Value = 1
$ws.Cells.Item(206, 10).Value = "H"
$ws.Cells.Item(206, 11).Value = 2.1
$ws.Cells.Item(206, 12).Value = 3.6
$ws.Cells.Item(206, 13).Value = 2.8
$ws.Cells.Item(206, 14).Value = 2.1
$ws.Cells.Item(206, 15).Value = 3.6
$ws.Cells.Item(206, 16).Value = 2.9
$ws.Cells.Item(206, 17).Value = -0.25
$ws.Cells.Item(206, 18).Value = 1.925
$ws.Cells.Item(206, 19).Value = 1.925
$ws.Cells.Item(206, 20).Value = 2.75
$ws.Cells.Item(206, 21).Value = 1.85
$ws.Cells.Item(206, 22).Value = 2
$ws.Cells.Item(206, 23).Value = 1.1
$ws.Cells.Item(206, 24).Value = -1
$ws.Cells.Item(206, 25).Value = -1
$ws.Cells.Item(206, 26).Value = 0.925
$ws.Cells.Item(206, 27).Value = -1
$ws.Cells.Item(206, 28).Value = 0.425
$ws.Cells.Item(206, 29).Value = -0.5

# row 214
$ws.Cells.Item(214, 2).Value = 6884507
$ws.Cells.Item(214, 5).Value = 45359.64583333334
$ws.Cells.Item(214, 6).Value = "Alemannia Aachen"
$ws.Cells.Item(214, 7).Value = "Fortuna Dusseldorf II"
$ws.Cells.Item(214, 8).Value = 2
$ws.Cells.Item(214, 9).Value = 1
$ws.Cells.Item(214, 10).Value = "H"
$ws.Cells.Item(214, 11).Value = 1.363
$ws.Cells.Item(214, 12).Value = 4.75
$ws.Cells.Item(214, 13).Value = 6
$ws.Cells.Item(214, 14).Value = 1.571
$ws.Cells.Item(214, 15).Value = 3.8
$ws.Cells.Item(214, 16).Value = 4.5
$ws.Cells.Item(214, 17).Value = -0.75
$ws.Cells.Item(214, 18).Value = 1.8
$ws.Cells.Item(214, 19).Value = 2
$ws.Cells.Item(214, 20).Value = 2.75
$ws.Cells.Item(214, 21).Value = 1.9
$ws.Cells.Item(214, 22).Value = 1.9
$ws.Cells.Item(214, 23).Value = 0.571
$ws.Cells.Item(214, 24).Value = -1
$ws.Cells.Item(214, 25).Value = -1
$ws.Cells.Item(214, 26).Value = 0.4
$ws.Cells.Item(214, 27).Value = -0.5
$ws.Cells.Item(214, 28).Value = 0.45
$ws.Cells.Item(214, 29).Value = -0.5

# row 215
$ws.Cells.Item(215, 2).Value = 6885380
$ws.Cells.Item(215, 5).Value = 45359.64583333334
$ws.Cells.Item(215, 6).Value = "SSVg Velbert"
$ws.Cells.Item(215, 7).Value = "Duren"
$ws.Cells.Item(215, 8).Value = 2
$ws.Cells.Item(215, 9).Value = 1
$ws.Cells.Item(215, 10).Value = "H"
$ws.Cells.Item(215, 11).Value = 3.25
$ws.Cells.Item(215, 12).Value = 4
$ws.Cells.Item(215, 13).Value = 1.8
$ws.Cells.Item(215, 14).Value = 3
$ws.Cells.Item(215, 15).Value = 3.6
$ws.Cells.Item(215, 16).Value = 1.95
$ws.Cells.Item(215, 17).Value = 0.25
$ws.Cells.Item(215, 18).Value = 2
$ws.Cells.Item(215, 19).Value = 1.8
$ws.Cells.Item(215, 20).Value = 3
$ws.Cells.Item(215, 21).Value = 1.9
$ws.Cells.Item(215, 22).Value = 1.9
$ws.Cells.Item(215, 23).Value = 2
$ws.Cells.Item(215, 24).Value = -1
$ws.Cells.Item(215, 25).Value = -1
$ws.Cells.Item(215, 26).Value = 1
$ws.Cells.Item(215, 27).Value = -1
$ws.Cells.Item(215, 28).Value = 0
$ws.Cells.Item(215, 29).Value = -0

# row 226
$ws.Cells.Item(226, 2).Value = 6886993
$ws.Cells.Item(226, 5).Value = 45367.41666666666
$ws.Cells.Item(226, 6).Value = "SC Fortuna Kln"
$ws.Cells.Item(226, 7).Value = "Fortuna Dusseldorf II"
$ws.Cells.Item(226, 8).Value = 4
$ws.Cells.Item(226, 9).Value = 1
$ws.Cells.Item(226, 10).Value = "H"
$ws.Cells.Item(226, 11).Value = 1.727
$ws.Cells.Item(226, 12).Value = 3.75
$ws.Cells.Item(226, 13).Value = 3.75
$ws.Cells.Item(226, 14).Value = 1.533
$ws.Cells.Item(226, 15).Value = 4
$ws.Cells.Item(226, 16).Value = 4.75
$ws.Cells.Item(226, 17).Value = -1
$ws.Cells.Item(226, 18).Value = 1.925
$ws.Cells.Item(226, 19).Value = 1.875
$ws.Cells.Item(226, 20).Value = 3
$ws.Cells.Item(226, 21).Value = 1.925
$ws.Cells.Item(226, 22).Value = 1.875
$ws.Cells.Item(226, 23).Value = 0.5329999999999999
$ws.Cells.Item(226, 24).Value = -1
$ws.Cells.Item(226, 25).Value = -1
$ws.Cells.Item(226, 26).Value = 0.925
$ws.Cells.Item(226, 27).Value = -1
$ws.Cells.Item(226, 28).Value = 0.925
$ws.Cells.Item(226, 29).Value = -1

# row 227
$ws.Cells.Item(227, 2).Value = 6884542
$ws.Cells.Item(227, 5).Value = 45367.41666666666
$ws.Cells.Item(227, 6).Value = "Wuppertaler"
$ws.Cells.Item(227, 7).Value = "SSVg Velbert"
$ws.Cells.Item(227, 8).Value = 5
$ws.Cells.Item(227, 9).Value = 0
$ws.Cells.Item(227, 10).Value = "H"
$ws.Cells.Item(227, 11).Value = 1.333
$ws.Cells.Item(227, 12).Value = 5
$ws.Cells.Item(227, 13).Value = 6
$ws.Cells.Item(227, 14).Value = 1.333
$ws.Cells.Item(227, 15).Value = 5.5
$ws.Cells.Item(227, 16).Value = 5.75
$ws.Cells.Item(227, 17).Value = -1.5
$ws.Cells.Item(227, 18).Value = 1.925
$ws.Cells.Item(227, 19).Value = 1.875
$ws.Cells.Item(227, 20).Value = 3.25
$ws.Cells.Item(227, 21).Value = 1.85
$ws.Cells.Item(227, 22).Value = 1.95
$ws.Cells.Item(227, 23).Value = 0.333
$ws.Cells.Item(227, 24).Value = -1
$ws.Cells.Item(227, 25).Value = -1
$ws.Cells.Item(227, 26).Value = 0.925
$ws.Cells.Item(227, 27).Value = -1
$ws.Cells.Item(227, 28).Value = 0.8500000000000001
$ws.Cells.Item(227, 29).Value = -1

# row 229
$ws.Cells.Item(229, 2).Value = 6886992
$ws.Cells.Item(229, 5).Value = 45367.41666666666
$ws.Cells.Item(229, 6).Value = "SV Rodinghausen"
$ws.Cells.Item(229, 7).Value = "FC WegbergBeeck"
$ws.Cells.Item(229, 8).Value = 5
$ws.Cells.Item(229, 9).Value = 1
$ws.Cells.Item(229, 10).Value = "H"
$ws.Cells.Item(229, 11).Value = 1.444
$ws.Cells.Item(229, 12).Value = 4.333
$ws.Cells.Item(229, 13).Value = 5.25
$ws.Cells.Item(229, 14).Value = 1.3
$ws.Cells.Item(229, 15).Value = 4.75
$ws.Cells.Item(229, 16).Value = 7.5
$ws.Cells.Item(229, 17).Value = -1.5
$ws.Cells.Item(229, 18).Value = 1.925
$ws.Cells.Item(229, 19).Value = 1.875
$ws.Cells.Item(229, 20).Value = 3
$ws.Cells.Item(229, 21).Value = 1.95
$ws.Cells.Item(229, 22).Value = 1.85
$ws.Cells.Item(229, 23).Value = 0.3
$ws.Cells.Item(229, 24).Value = -1
$ws.Cells.Item(229, 25).Value = -1
$ws.Cells.Item(229, 26).Value = 0.925
$ws.Cells.Item(229, 27).Value = -1
$ws.Cells.Item(229, 28).Value = 0.95
$ws.Cells.Item(229, 29).Value = -1

# row 230
$ws.Cells.Item(230, 2).Value = 6886994
$ws.Cells.Item(230, 5).Value = 45367.41666666666
$ws.Cells.Item(230, 6).Value = "Duren"
$ws.Cells.Item(230, 7).Value = "SV Lippstadt 08"
$ws.Cells.Item(230, 8).Value = 2
$ws.Cells.Item(230, 9).Value = 0
$ws.Cells.Item(230, 10).Value = "H"
$ws.Cells.Item(230, 11).Value = 1.727
$ws.Cells.Item(230, 12).Value = 4
$ws.Cells.Item(230, 13).Value = 3.5
$ws.Cells.Item(230, 14).Value = 1.6
$ws.Cells.Item(230, 15).Value = 4
$ws.Cells.Item(230, 16).Value = 4.2
$ws.Cells.Item(230, 17).Value = -0.75
$ws.Cells.Item(230, 18).Value = 1.775
$ws.Cells.Item(230, 19).Value = 2.025
$ws.Cells.Item(230, 20).Value = 2.75
$ws.Cells.Item(230, 21).Value = 1.8
$ws.Cells.Item(230, 22).Value = 2
$ws.Cells.Item(230, 23).Value = 0.6000000000000001
$ws.Cells.Item(230, 24).Value = -1
$ws.Cells.Item(230, 25).Value = -1
$ws.Cells.Item(230, 26).Value = 0.7749999999999999
$ws.Cells.Item(230, 27).Value = -1
$ws.Cells.Item(230, 28).Value = -1
$ws.Cells.Item(230, 29).Value = 1

# row 231
$ws.Cells.Item(231, 2).Value = 6880549
$ws.Cells.Item(231, 5).Value = 45367.41666666666
$ws.Cells.Item(231, 6).Value = "1 FC Bocholt"
$ws.Cells.Item(231, 7).Value = "RotWeiss Oberhausen"
$ws.Cells.Item(231, 8).Value = 0
$ws.Cells.Item(231, 9).Value = 0
$ws.Cells.Item(231, 10).Value = "D"
$ws.Cells.Item(231, 11).Value = 2.1
$ws.Cells.Item(231, 12).Value = 3.6
$ws.Cells.Item(231, 13).Value = 2.8
$ws.Cells.Item(231, 14).Value = 3
$ws.Cells.Item(231, 15).Value = 3.4
$ws.Cells.Item(231, 16).Value = 2.2
$ws.Cells.Item(231, 17).Value = 0.25
$ws.Cells.Item(231, 18).Value = 1.875
$ws.Cells.Item(231, 19).Value = 1.975
$ws.Cells.Item(231, 20).Value = 2.5
$ws.Cells.Item(231, 21).Value = 1.825
$ws.Cells.Item(231, 22).Value = 2.025
$ws.Cells.Item(231, 23).Value = -1
$ws.Cells.Item(231, 24).Value = 2.4
$ws.Cells.Item(231, 25).Value = -1
$ws.Cells.Item(231, 26).Value = 0.4375
$ws.Cells.Item(231, 27).Value = -0.5
$ws.Cells.Item(231, 28).Value = -1
$ws.Cells.Item(231, 29).Value = 1.025

# row 238
$ws.Cells.Item(238, 2).Value = 6884509
$ws.Cells.Item(238, 5).Value = 45381.41666666666
$ws.Cells.Item(238, 6).Value = "Alemannia Aachen"
$ws.Cells.Item(238, 7).Value = "SC Fortuna Kln"
$ws.Cells.Item(238, 8).Value = 1
$ws.Cells.Item(238, 9).Value = 0
$ws.Cells.Item(238, 10).Value = "H"
$ws.Cells.Item(238, 11).Value = 1.833
$ws.Cells.Item(238, 12).Value = 3.75
$ws.Cells.Item(238, 13).Value = 3.3
$ws.Cells.Item(238, 14).Value = 1.7
$ws.Cells.Item(238, 15).Value = 3.8
$ws.Cells.Item(238, 16).Value = 4
$ws.Cells.Item(238, 17).Value = -0.75
$ws.Cells.Item(238, 18).Value = 1.925
$ws.Cells.Item(238, 19).Value = 1.875
$ws.Cells.Item(238, 20).Value = 2.5
$ws.Cells.Item(238, 21).Value = 1.825
$ws.Cells.Item(238, 22).Value = 1.975
$ws.Cells.Item(238, 23).Value = 0.7
$ws.Cells.Item(238, 24).Value = -1
$ws.Cells.Item(238, 25).Value = -1
$ws.Cells.Item(238, 26).Value = 0.4625
$ws.Cells.Item(238, 27).Value = -0.5
$ws.Cells.Item(238, 28).Value = -1
$ws.Cells.Item(238, 29).Value = 0.9750000000000001

# row 239
$ws.Cells.Item(239, 2).Value = 6882756
$ws.Cells.Item(239, 5).Value = 45381.41666666666
$ws.Cells.Item(239, 6).Value = "Gutersloh 2000"
$ws.Cells.Item(239, 7).Value = "1 FC Bocholt"
$ws.Cells.Item(239, 8).Value = 1
$ws.Cells.Item(239, 9).Value = 0
$ws.Cells.Item(239, 10).Value = "H"
$ws.Cells.Item(239, 11).Value = 3.1
$ws.Cells.Item(239, 12).Value = 3.75
$ws.Cells.Item(239, 13).Value = 1.909
$ws.Cells.Item(239, 14).Value = 2.75
$ws.Cells.Item(239, 15).Value = 3.3
$ws.Cells.Item(239, 16).Value = 2.25
$ws.Cells.Item(239, 17).Value = 0.25
$ws.Cells.Item(239, 18).Value = 1.775
$ws.Cells.Item(239, 19).Value = 2.025
$ws.Cells.Item(239, 20).Value = 2.5
$ws.Cells.Item(239, 21).Value = 1.975
$ws.Cells.Item(239, 22).Value = 1.825
$ws.Cells.Item(239, 23).Value = 1.75
$ws.Cells.Item(239, 24).Value = -1
$ws.Cells.Item(239, 25).Value = -1
$ws.Cells.Item(239, 26).Value = 0.7749999999999999
$ws.Cells.Item(239, 27).Value = -1
$ws.Cells.Item(239, 28).Value = -1
$ws.Cells.Item(239, 29).Value = 0.825

# row 240
$ws.Cells.Item(240, 2).Value = 6884543
$ws.Cells.Item(240, 5).Value = 45381.41666666666
$ws.Cells.Item(240, 6).Value = "SV Lippstadt 08"
$ws.Cells.Item(240, 7).Value = "Wuppertaler"
$ws.Cells.Item(240, 8).Value = 0
$ws.Cells.Item(240, 9).Value = 1
$ws.Cells.Item(240, 10).Value = "A"
$ws.Cells.Item(240, 11).Value = 3.1
$ws.Cells.Item(240, 12).Value = 3.75
$ws.Cells.Item(240, 13).Value = 1.909
$ws.Cells.Item(240, 14).Value = 3
$ws.Cells.Item(240, 15).Value = 3.75
$ws.Cells.Item(240, 16).Value = 1.95
$ws.Cells.Item(240, 17).Value = 0.5
$ws.Cells.Item(240, 18).Value = 1.825
$ws.Cells.Item(240, 19).Value = 2.025
$ws.Cells.Item(240, 20).Value = 3
$ws.Cells.Item(240, 21).Value = 1.925
$ws.Cells.Item(240, 22).Value = 1.925
$ws.Cells.Item(240, 23).Value = -1
$ws.Cells.Item(240, 24).Value = -1
$ws.Cells.Item(240, 25).Value = 0.95
$ws.Cells.Item(240, 26).Value = -1
$ws.Cells.Item(240, 27).Value = 1.025
$ws.Cells.Item(240, 28).Value = -1
$ws.Cells.Item(240, 29).Value = 0.925

# row 241
$ws.Cells.Item(241, 2).Value = 6886996
$ws.Cells.Item(241, 5).Value = 45382.375
$ws.Cells.Item(241, 6).Value = "FC WegbergBeeck"
$ws.Cells.Item(241, 7).Value = "Cologne II"
$ws.Cells.Item(241, 8).Value = 0
$ws.Cells.Item(241, 9).Value = 1
$ws.Cells.Item(241, 10).Value = "A"
$ws.Cells.Item(241, 11).Value = 3
$ws.Cells.Item(241, 12).Value = 3.6
$ws.Cells.Item(241, 13).Value = 2
$ws.Cells.Item(241, 14).Value = 3.2
$ws.Cells.Item(241, 15).Value = 3.75
$ws.Cells.Item(241, 16).Value = 1.909
$ws.Cells.Item(241, 17).Value = 0.5
$ws.Cells.Item(241, 18).Value = 1.825
$ws.Cells.Item(241, 19).Value = 1.975
$ws.Cells.Item(241, 20).Value = 3
$ws.Cells.Item(241, 21).Value = 1.95
$ws.Cells.Item(241, 22).Value = 1.85
$ws.Cells.Item(241, 23).Value = -1
$ws.Cells.Item(241, 24).Value = -1
$ws.Cells.Item(241, 25).Value = 0.909
$ws.Cells.Item(241, 26).Value = -1
$ws.Cells.Item(241, 27).Value = 0.9750000000000001
$ws.Cells.Item(241, 28).Value = -1
$ws.Cells.Item(241, 29).Value = 0.8500000000000001

# row 242
$ws.Cells.Item(242, 2).Value = 6886997
$ws.Cells.Item(242, 5).Value = 45382.375
$ws.Cells.Item(242, 6).Value = "SC Paderborn 07 II"
$ws.Cells.Item(242, 7).Value = "SV Rodinghausen"
$ws.Cells.Item(242, 8).Value = 0
$ws.Cells.Item(242, 9).Value = 0
$ws.Cells.Item(242, 10).Value = "D"
$ws.Cells.Item(242, 11).Value = 2.7
$ws.Cells.Item(242, 12).Value = 3.6
$ws.Cells.Item(242, 13).Value = 2.15
$ws.Cells.Item(242, 14).Value = 3.1
$ws.Cells.Item(242, 15).Value = 3.5
$ws.Cells.Item(242, 16).Value = 1.95
$ws.Cells.Item(242, 17).Value = 0.5
$ws.Cells.Item(242, 18).Value = 1.775
$ws.Cells.Item(242, 19).Value = 2.025
$ws.Cells.Item(242, 20).Value = 2.5
$ws.Cells.Item(242, 21).Value = 1.8
$ws.Cells.Item(242, 22).Value = 2
$ws.Cells.Item(242, 23).Value = -1
$ws.Cells.Item(242, 24).Value = 2.5
$ws.Cells.Item(242, 25).Value = -1
$ws.Cells.Item(242, 26).Value = 0.7749999999999999
$ws.Cells.Item(242, 27).Value = -1
$ws.Cells.Item(242, 28).Value = -1
$ws.Cells.Item(242, 29).Value = 1

# --- Section 2: append two brand-new fixtures (rows 243-244) ---
# Copy the A/E number formats (bold-bordered id cell / date cell) from the last row
# down into the two new rows before filling in values, so styling matches the rest of the table.
$ws.Range("A242").Copy()
$ws.Range("A243:A244").PasteSpecial(-4122)
$ws.Range("E242").Copy()
$ws.Range("E243:E244").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# row 243
$ws.Cells.Item(243, 1).Value = 241
$ws.Cells.Item(243, 2).Value = 6887000
$ws.Cells.Item(243, 3).Value = "Germany Regionalliga West"
$ws.Cells.Item(243, 4).Value = "Germany Regionalliga West"
$ws.Cells.Item(243, 5).Value = 45387.60416666666
$ws.Cells.Item(243, 6).Value = "SV Rodinghausen"
$ws.Cells.Item(243, 7).Value = "SC Wiedenbruck"
$ws.Cells.Item(243, 11).Value = 2.05
$ws.Cells.Item(243, 12).Value = 3.4
$ws.Cells.Item(243, 13).Value = 3
$ws.Cells.Item(243, 14).Value = 1.666
$ws.Cells.Item(243, 15).Value = 3.6
$ws.Cells.Item(243, 16).Value = 4.2
$ws.Cells.Item(243, 17).Value = -0.5
$ws.Cells.Item(243, 18).Value = 1.725
$ws.Cells.Item(243, 19).Value = 2.15
$ws.Cells.Item(243, 20).Value = 2.5
$ws.Cells.Item(243, 21).Value = 1.85
$ws.Cells.Item(243, 22).Value = 2
$ws.Cells.Item(243, 23).Value = 0
$ws.Cells.Item(243, 24).Value = 0
$ws.Cells.Item(243, 25).Value = 0
$ws.Cells.Item(243, 26).Value = 0
$ws.Cells.Item(243, 27).Value = 0

# row 244
$ws.Cells.Item(244, 1).Value = 242
$ws.Cells.Item(244, 2).Value = 6887001
$ws.Cells.Item(244, 3).Value = "Germany Regionalliga West"
$ws.Cells.Item(244, 4).Value = "Germany Regionalliga West"
$ws.Cells.Item(244, 5).Value = 45387.60416666666
$ws.Cells.Item(244, 6).Value = "Cologne II"
$ws.Cells.Item(244, 7).Value = "SC Paderborn 07 II"
$ws.Cells.Item(244, 11).Value = 1.95
$ws.Cells.Item(244, 12).Value = 3.5
$ws.Cells.Item(244, 13).Value = 3.2
$ws.Cells.Item(244, 14).Value = 1.909
$ws.Cells.Item(244, 15).Value = 3.5
$ws.Cells.Item(244, 16).Value = 3.4
$ws.Cells.Item(244, 17).Value = -0.5
$ws.Cells.Item(244, 18).Value = 1.975
$ws.Cells.Item(244, 19).Value = 1.875
$ws.Cells.Item(244, 20).Value = 2.75
$ws.Cells.Item(244, 21).Value = 1.8
$ws.Cells.Item(244, 22).Value = 2.05
$ws.Cells.Item(244, 23).Value = 0
$ws.Cells.Item(244, 24).Value = 0
$ws.Cells.Item(244, 25).Value = 0
$ws.Cells.Item(244, 26).Value = 0
$ws.Cells.Item(244, 27).Value = 0

